$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Write PASS into the Result column for the rows that validated successfully
$ws.Range("S2").Value = "PASS"
$ws.Range("S3").Value = "PASS"
$ws.Range("S5").Value = "PASS"
$ws.Range("S6").Value = "PASS"
$ws.Range("S8").Value = "PASS"

# The autofilter now only covers through column R
$ws.AutoFilterMode = $false
$ws.Range("A1:R8").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=TestData!`$A`$1:`$R`$8"
    }
}

# Selection left where the user was working after the edit
$ws.Range("S2:S10").Select()
